$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen the new date column F (closest achievable pixel-snapped width to 19.45)
$ws.Columns.Item(6).ColumnWidth = 18.6

# Re-type F2:F11 from plain "2024" numbers to a real date serial (2024-06-24),
# formatted with the built-in short-date numeric format (numFmtId 14).
$ws.Range("F2").Value = 45467
$ws.Range("F2").NumberFormat = "mm-dd-yy"

# Share that exact same style across the rest of the column instead of
# re-deriving a fresh style per cell (keeps a single shared cellXfs entry,
# matching how a uniform style column really looks).
$ws.Range("F2").Copy()
$ws.Range("F3:F11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("F3").Value = 45467
$ws.Range("F4").Value = 45467
$ws.Range("F5").Value = 45467
$ws.Range("F6").Value = 45467
$ws.Range("F7").Value = 45467
$ws.Range("F8").Value = 45467
$ws.Range("F9").Value = 45467
$ws.Range("F10").Value = 45467
$ws.Range("F11").Value = 45467

# Row heights settle back down now the dates are short, fixed-width values.
$ws.Rows.Item(5).RowHeight = 31
$ws.Rows.Item(6).RowHeight = 15.5
$ws.Rows.Item(7).RowHeight = 31
$ws.Rows.Item(8).RowHeight = 15.5
$ws.Rows.Item(9).RowHeight = 31
$ws.Rows.Item(10).RowHeight = 15.5
$ws.Rows.Item(11).RowHeight = 31

$ws.Range("E27").Select()
